$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.026976621435872
$ws.Range("D2").Value = 1.02956058023653
$ws.Range("E2").Value = 1.036641762943746
$ws.Range("F2").Value = 1.046593706844036
$ws.Range("I2").Value = 1.029305136208412
$ws.Range("J2").Value = 1.03213706313554
$ws.Range("K2").Value = 1.032374135958401
$ws.Range("L2").Value = 1.039434922222618
$ws.Range("M2").Value = 1.049358707036724
$ws.Range("N2").Value = 1.033602816152981
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.027957525184079
$ws.Range("D3").Value = 1.03024396543975
$ws.Range("E3").Value = 1.037586974853386
$ws.Range("F3").Value = 1.047791776093722
$ws.Range("I3").Value = 1.029422768792544
$ws.Range("J3").Value = 1.032757794265601
$ws.Range("K3").Value = 1.032866028154563
$ws.Range("L3").Value = 1.040189418138664
$ws.Range("M3").Value = 1.050367446765378
$ws.Range("N3").Value = 1.034224428792434
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.02859251344822
$ws.Range("D4").Value = 1.030686121826785
$ws.Range("E4").Value = 1.038199261807337
$ws.Range("F4").Value = 1.048568167141819
$ws.Range("I4").Value = 1.029497481391321
$ws.Range("J4").Value = 1.033159143959063
$ws.Range("K4").Value = 1.033183599474197
$ws.Range("L4").Value = 1.040677677490814
$ws.Range("M4").Value = 1.051020758284741
$ws.Range("N4").Value = 1.034626348448495
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.028859528596911
$ws.Range("D5").Value = 1.030871993750427
$ws.Range("E5").Value = 1.038456827348609
$ws.Range("F5").Value = 1.048894839413128
$ws.Range("I5").Value = 1.029528554110763
$ws.Range("J5").Value = 1.033327797766417
$ws.Range("K5").Value = 1.033316934306799
$ws.Range("L5").Value = 1.040882953156755
$ws.Range("M5").Value = 1.051295551200866
$ws.Range("N5").Value = 1.034795241763601
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.028904365477165
$ws.Range("D6").Value = 1.030903201807772
$ws.Range("E6").Value = 1.038500083105104
$ws.Range("F6").Value = 1.04894970538866
$ws.Range("I6").Value = 1.029533751614158
$ws.Range("J6").Value = 1.0333561111486
$ws.Range("K6").Value = 1.033339311694222
$ws.Range("L6").Value = 1.04091742049902
$ws.Range("M6").Value = 1.05134169838398
$ws.Range("N6").Value = 1.034823595354034
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.028596081058613
$ws.Range("D7").Value = 1.030688605499721
$ws.Range("E7").Value = 1.038202702782457
$ws.Range("F7").Value = 1.048572531062728
$ws.Range("I7").Value = 1.029497897909122
$ws.Range("J7").Value = 1.033161397809483
$ws.Range("K7").Value = 1.033185381778173
$ws.Range("L7").Value = 1.040680420350062
$ws.Range("M7").Value = 1.05102442952852
$ws.Range("N7").Value = 1.034628605499641
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.027308064622169
$ws.Range("D8").Value = 1.029791540957082
$ws.Range("E8").Value = 1.036961062251733
$ws.Range("F8").Value = 1.046998359583137
$ws.Range("I8").Value = 1.029345181037698
$ws.Range("J8").Value = 1.032346904875704
$ws.Range("K8").Value = 1.032540520975797
$ws.Range("L8").Value = 1.03968989695907
$ws.Range("M8").Value = 1.049699493077365
$ws.Range("N8").Value = 1.033812955892473
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.02504056388426
$ws.Range("D9").Value = 1.028210546676084
$ws.Range("E9").Value = 1.034778307136867
$ws.Range("F9").Value = 1.044233364911079
$ws.Range("I9").Value = 1.029065344461944
$ws.Range("J9").Value = 1.030909356681318
$ws.Range("K9").Value = 1.031398745801676
$ws.Range("L9").Value = 1.037944875390918
$ws.Range("M9").Value = 1.047369316259821
$ws.Range("N9").Value = 1.032373366214762
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.023530366306495
$ws.Range("D10").Value = 1.027156448664628
$ws.Range("E10").Value = 1.033326654779589
$ws.Range("F10").Value = 1.042396029494541
$ws.Range("I10").Value = 1.028871598512176
$ws.Range("J10").Value = 1.029949475167546
$ws.Range("K10").Value = 1.030633944653139
$ws.Range("L10").Value = 1.036781832172202
$ws.Range("M10").Value = 1.045818930386779
$ws.Range("N10").Value = 1.031412121559139
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.022876787177276
$ws.Range("D11").Value = 1.026700002364288
$ws.Range("E11").Value = 1.032698915281215
$ws.Range("F11").Value = 1.041601866615615
$ws.Range("I11").Value = 1.028786004747578
$ws.Range("J11").Value = 1.029533483331833
$ws.Range("K11").Value = 1.030301928458542
$ws.Range("L11").Value = 1.036278299802555
$ws.Range("M11").Value = 1.045148326597921
$ws.Range("N11").Value = 1.030995538967305
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.022634071105623
$ws.Range("D12").Value = 1.026530456728948
$ws.Range("E12").Value = 1.032465870964337
$ws.Range("F12").Value = 1.04130709187564
$ws.Range("I12").Value = 1.028753956327508
$ws.Range("J12").Value = 1.02937891234629
$ws.Range("K12").Value = 1.030178475514376
$ws.Range("L12").Value = 1.036091276846211
$ws.Range("M12").Value = 1.04489934337771
$ws.Range("N12").Value = 1.030840748473233
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.022686132154949
$ws.Range("D13").Value = 1.026566824878218
$ws.Range("E13").Value = 1.032515854047014
$ws.Range("F13").Value = 1.041370312431358
$ws.Range("I13").Value = 1.028760842359193
$ws.Range("J13").Value = 1.029412070741882
$ws.Range("K13").Value = 1.030204962341522
$ws.Range("L13").Value = 1.036131393383302
$ws.Range("M13").Value = 1.044952746193331
$ws.Range("N13").Value = 1.03087395395755
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.022856723122199
$ws.Range("D14").Value = 1.02668598768651
$ws.Range("E14").Value = 1.032679649190036
$ws.Range("F14").Value = 1.041577496104034
$ws.Range("I14").Value = 1.028783360816548
$ws.Range("J14").Value = 1.02952070752613
$ws.Range("K14").Value = 1.030291726393481
$ws.Range("L14").Value = 1.036262840191969
$ws.Range("M14").Value = 1.045127743344951
$ws.Range("N14").Value = 1.030982745018493
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.02296183678827
$ws.Range("D15").Value = 1.026759407704728
$ws.Range("E15").Value = 1.032780585507299
$ws.Range("F15").Value = 1.04170517699098
$ws.Range("I15").Value = 1.028797201391919
$ws.Range("J15").Value = 1.029587635213959
$ws.Range("K15").Value = 1.030345167735589
$ws.Range("L15").Value = 1.036343830430238
$ws.Range("M15").Value = 1.045235579299693
$ws.Range("N15").Value = 1.031049767751315
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.023573749522306
$ws.Range("D16").Value = 1.027186741280798
$ws.Range("E16").Value = 1.033368333433759
$ws.Range("F16").Value = 1.042448765296281
$ws.Range("I16").Value = 1.028877243294739
$ws.Range("J16").Value = 1.029977075701694
$ws.Range("K16").Value = 1.030655961577868
$ws.Range("L16").Value = 1.03681525153105
$ws.Range("M16").Value = 1.045863451411118
$ws.Range("N16").Value = 1.031439761289212
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.023957679145908
$ws.Range("D17").Value = 1.027454793118096
$ws.Range("E17").Value = 1.033737236269084
$ws.Range("F17").Value = 1.042915577595346
$ws.Range("I17").Value = 1.028926996474743
$ws.Range("J17").Value = 1.030221266263749
$ws.Range("K17").Value = 1.030850686663015
$ws.Range("L17").Value = 1.037110981374608
$ws.Range("M17").Value = 1.046257492835367
$ws.Range("N17").Value = 1.031684298629874
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.024181652202407
$ws.Range("D18").Value = 1.027611141709059
$ws.Range("E18").Value = 1.033952491633172
$ws.Range("F18").Value = 1.043187997881054
$ws.Range("I18").Value = 1.028955852563795
$ws.Range("J18").Value = 1.030363664026334
$ws.Range("K18").Value = 1.030964184234344
$ws.Range("L18").Value = 1.037283482606301
$ws.Range("M18").Value = 1.046487400328558
$ws.Range("N18").Value = 1.031826898613616
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.024258026883995
$ws.Range("D19").Value = 1.027664452248939
$ws.Range("E19").Value = 1.03402590180094
$ws.Range("F19").Value = 1.043280909321528
$ws.Range("I19").Value = 1.028965663898403
$ws.Range("J19").Value = 1.030412212122813
$ws.Range("K19").Value = 1.031002869996494
$ws.Range("L19").Value = 1.037342302246134
$ws.Range("M19").Value = 1.046565804675659
$ws.Range("N19").Value = 1.031875515653961
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.023916483673694
$ws.Range("D20").Value = 1.027426033841158
$ws.Range("E20").Value = 1.033697648147136
$ws.Range("F20").Value = 1.042865478927009
$ws.Range("I20").Value = 1.02892167540006
$ws.Range("J20").Value = 1.030195070488147
$ws.Range("K20").Value = 1.030829803009407
$ws.Range("L20").Value = 1.037079251632257
$ws.Range("M20").Value = 1.046215208698855
$ws.Range("N20").Value = 1.031658065653266
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.022806486887749
$ws.Range("D21").Value = 1.026650897232752
$ws.Range("E21").Value = 1.032631412117964
$ws.Range("F21").Value = 1.041516479803961
$ws.Range("I21").Value = 1.028776736730535
$ws.Range("J21").Value = 1.029488718162469
$ws.Range("K21").Value = 1.030266180052528
$ws.Range("L21").Value = 1.036224132065401
$ws.Range("M21").Value = 1.045076208038212
$ws.Range("N21").Value = 1.030950710226266
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.022108890686034
$ws.Range("D22").Value = 1.026163531927309
$ws.Range("E22").Value = 1.031961756958803
$ws.Range("F22").Value = 1.040669541238313
$ws.Range("I22").Value = 1.02868413238975
$ws.Range("J22").Value = 1.029044299478179
$ws.Range("K22").Value = 1.029911071485355
$ws.Range("L22").Value = 1.035686550534523
$ws.Range("M22").Value = 1.044360703189233
$ws.Range("N22").Value = 1.030505660416468
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.022478670832985
$ws.Range("D23").Value = 1.026421893719385
$ws.Range("E23").Value = 1.032316684472549
$ws.Range("F23").Value = 1.041118402718801
$ws.Range("I23").Value = 1.028733363445371
$ws.Range("J23").Value = 1.029279923155323
$ws.Range("K23").Value = 1.030099390830363
$ws.Range("L23").Value = 1.03597152627253
$ws.Range("M23").Value = 1.04473994606114
$ws.Range("N23").Value = 1.030741618706261
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.023935098031062
$ws.Range("D24").Value = 1.027439028925225
$ws.Range("E24").Value = 1.033715536066387
$ws.Range("F24").Value = 1.042888115936649
$ws.Range("I24").Value = 1.028924080271727
$ws.Range("J24").Value = 1.030206907338026
$ws.Range("K24").Value = 1.030839239687793
$ws.Range("L24").Value = 1.037093588915316
$ws.Range("M24").Value = 1.046234314864286
$ws.Range("N24").Value = 1.03166991931283
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.025626509985362
$ws.Range("D25").Value = 1.028619294477664
$ws.Range("E25").Value = 1.035341984488402
$ws.Range("F25").Value = 1.04494712671842
$ws.Range("I25").Value = 1.029138957801592
$ws.Range("J25").Value = 1.031281267158212
$ws.Range("K25").Value = 1.031694563006964
$ws.Range("L25").Value = 1.03839595388061
$ws.Range("M25").Value = 1.047971183961924
$ws.Range("N25").Value = 1.032745804847192
